# start work on get-hplldpneighbor
#
# The "Vlan" class hierarchy table is extended with a new "Interface" class
# block (rows 10-22), and the property type for Vlan.DhcpRelayList is
# updated from "array[string]" to "list<string>". Work-in-progress rows
# 23/24 start sketching the next class's properties (values only, no types
# yet), matching the "start work on ..." commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Vlan class: fix up DhcpRelayEnabled/DhcpRelayList rows --------------
$ws.Range("B7").Value = "DhcpRelayEnabled"
$ws.Range("C7").Value = "bool"
$ws.Range("B8").Value = "DhcpRelayList"
$ws.Range("C8").Value = "list<string>"

# --- New "Interface" class block -----------------------------------------
$ws.Range("A10").Value = "Interface"
$ws.Range("B10").Value = "class"

$ws.Range("B11").Value = "Name"
$ws.Range("C11").Value = "string"

$ws.Range("B12").Value = "Description"
$ws.Range("C12").Value = "string"

$ws.Range("B13").Value = "IpAddress"
$ws.Range("C13").Value = "string"

$ws.Range("B14").Value = "PimSmEnabled"
$ws.Range("C14").Value = "bool"

$ws.Range("B15").Value = "PortLinkType"
$ws.Range("C15").Value = "string"

$ws.Range("B16").Value = "PortLinkMode"
$ws.Range("C16").Value = "string"

$ws.Range("B17").Value = "LinkAggMode"
$ws.Range("C17").Value = "string"

$ws.Range("B18").Value = "LinkAggGroup"
$ws.Range("C18").Value = "int"

$ws.Range("B19").Value = "MadEnabled"
$ws.Range("C19").Value = "bool"

$ws.Range("B20").Value = "PermittedVlans"
$ws.Range("C20").Value = "list<int>"

$ws.Range("B21").Value = "Pvid"
$ws.Range("C21").Value = "int"

$ws.Range("B22").Value = "IsShutdown"
$ws.Range("C22").Value = "bool"

# --- Start of next block (work in progress, no types filled in yet) ------
$ws.Range("B23").Value = "DhcpRelayEnabled"
$ws.Range("B24").Value = "DhcpRelayList"

# Leave the selection where the author was last working.
$ws.Range("B18").Select()
